{"js": "// Re-anchor the three heading bookmarks so each one spans the full heading\n// text (instead of being an empty, zero-length bookmark sitting just before\n// the text), and rename the \"Acknowledgement by Authorized Investigator\"\n// bookmark to its regenerated slug.\n\nasync function rewrapBookmark(oldName, newName) {\n  // Locate the (currently empty) bookmark and the paragraph it sits in.\n  const bm = context.document.bookmarks.getByName(oldName);\n  const bmRange = bm.getRange();\n  const paragraphs = bmRange.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const heading = paragraphs.items[0];\n  // \"Content\" excludes the paragraph mark, so the new bookmark end stays\n  // inside this paragraph instead of spilling into the next one.\n  const target = heading.getRange(Word.RangeLocation.content);\n\n  // Remove the old (empty) bookmark, then re-insert one that wraps the\n  // heading's full text under the (possibly new) name.\n  context.document.deleteBookmark(oldName);\n  await context.sync();\n\n  target.insertBookmark(newName);\n  await context.sync();\n}\n\nawait rewrapBookmark(\"databrary-access-agreement\", \"databrary-access-agreement\");\nawait rewrapBookmark(\n  \"acknowledgement-by-authorized-investigator\",\n  \"X2fdad48066937982775bea765b740eea6efd90c\"\n);\nawait rewrapBookmark(\"approval-by-institution\", \"approval-by-institution\");\n", "ps1": "# Re-anchor the three heading bookmarks so each one spans the full heading\n# text (instead of being an empty, zero-length bookmark sitting just before\n# the text), and rename the \"Acknowledgement by Authorized Investigator\"\n# bookmark to its regenerated slug.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphAt([int]$pos) {\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Start -eq $pos) {\n            return $p\n        }\n    }\n    return $d.Paragraphs.Item(1)\n}\n\nfunction Rewrap-Bookmark([string]$oldName, [string]$newName) {\n    $bm = $d.Bookmarks.Item($oldName)\n    $bmStart = $bm.Range.Start\n\n    $para = Get-ParagraphAt($bmStart)\n    $target = $para.Range\n\n    # Trim a trailing paragraph mark, if the range object includes one, so\n    # the bookmark stays inside this paragraph instead of spilling into the\n    # next one.\n    $text = $target.Text\n    if ($text.Length -gt 0 -and $text.Substring($text.Length - 1) -eq [char]13) {\n        $target.MoveEnd(1, -1)\n    }\n\n    $bm.Delete()\n    $d.Bookmarks.Add($newName, $target)\n}\n\nRewrap-Bookmark \"databrary-access-agreement\" \"databrary-access-agreement\"\nRewrap-Bookmark \"acknowledgement-by-authorized-investigator\" \"X2fdad48066937982775bea765b740eea6efd90c\"\nRewrap-Bookmark \"approval-by-institution\" \"approval-by-institution\"\n"}
